$wb = $excel.ActiveWorkbook

# addListItem: update the Surat value (auto propagates via the C2 formula =A2)
$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "SuratP"

# createUser: bump the test user id (auto propagates via formulas in B2/F2)
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 62
$wsCreateUser.Range("D11").Select() | Out-Null

# Make addListItem the active/selected sheet & tab
$wsAddListItem.Activate()
